$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.446.60'
$ws.Range("E2").Value = '  +2.33%  '
$ws.Range("D3").Value = '2.041.84'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.85'
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.43'
$ws.Range("E7").Value = '  -1.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.394'
$ws.Range("E9").Value = '  +4.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0807'
$ws.Range("E10").Value = '  +2.59%  '
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.19'
$ws.Range("E12").Value = '  +6.95%  '
$ws.Range("D13").Value = '2.337.40'
$ws.Range("E13").Value = '  +3.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.851'
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.10'
$ws.Range("E15").Value = '  +2.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '2.043.08'
$ws.Range("E17").Value = '  +4.02%  '
$ws.Range("D18").Value = '37.340.54'
$ws.Range("E18").Value = '  +2.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.45'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '0.0₃0863'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.29'
$ws.Range("E21").Value = '  +4.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.00'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("E24").Value = '  +4.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.36'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.55'
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("E28").Value = '  -3.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.94'
$ws.Range("E29").Value = '  +3.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.38'
$ws.Range("E30").Value = '  +3.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0680'
$ws.Range("E32").Value = '  +10.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.81'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.53'
$ws.Range("E34").Value = '  +11.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.53'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.60'
$ws.Range("E36").Value = '  +5.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.47'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0983'
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0217'
$ws.Range("E42").Value = '  +3.65%  '
$ws.Range("E43").Value = '  +1.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.71'
$ws.Range("E44").Value = '  +4.92%  '
$ws.Range("D45").Value = '1.404.36'
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.83'
$ws.Range("E46").Value = '  +3.32%  '
$ws.Range("E47").Value = '  +3.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.52'
$ws.Range("E48").Value = '  +3.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.12'
$ws.Range("E49").Value = '  +14.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("E51").Value = '  +3.65%  '
